{"js": "const replacements = [\n  [\"2024-03-25 Monday\", \"2024-03-26 Tuesday\"],\n  [\"861\u00f79=\", \"513\u00f76=\"],\n  [\"645\u00f79=\", \"845\u00f75=\"],\n  [\"755\u00f76=\", \"168\u00f78=\"],\n  [\"143\u00f73=\", \"122\u00f76=\"],\n  [\"609\u00f79=\", \"201\u00f73=\"],\n  [\"724\u00f78=\", \"543\u00f73=\"],\n  [\"907\u00f73=\", \"974\u00f77=\"],\n  [\"571\u00f79=\", \"410\u00f72=\"],\n  [\"235\u00f76=\", \"473\u00f78=\"],\n  [\"684\u00f76=\", \"929\u00f73=\"],\n  [\"449\u00f74=\", \"708\u00f77=\"],\n  [\"598\u00f75=\", \"790\u00f76=\"],\n  [\"758\u00f77=\", \"360\u00f76=\"],\n  [\"109\u00f74=\", \"855\u00f77=\"],\n  [\"361\u00f76=\", \"958\u00f79=\"],\n  [\"276\u00f74=\", \"143\u00f76=\"],\n  [\"986\u00f76=\", \"565\u00f74=\"],\n  [\"696\u00f78=\", \"678\u00f77=\"],\n  [\"623\u00f72=\", \"812\u00f75=\"],\n  [\"254\u00f77=\", \"193\u00f72=\"],\n  [\"740\u00f79=\", \"924\u00f76=\"],\n  [\"500\u00f74=\", \"888\u00f72=\"],\n  [\"950\u00f78=\", \"881\u00f72=\"],\n  [\"838\u00f76=\", \"875\u00f77=\"],\n  [\"865\u00f77=\", \"740\u00f74=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-03-25 Monday\", \"2024-03-26 Tuesday\"),\n    @(\"861\u00f79=\", \"513\u00f76=\"),\n    @(\"645\u00f79=\", \"845\u00f75=\"),\n    @(\"755\u00f76=\", \"168\u00f78=\"),\n    @(\"143\u00f73=\", \"122\u00f76=\"),\n    @(\"609\u00f79=\", \"201\u00f73=\"),\n    @(\"724\u00f78=\", \"543\u00f73=\"),\n    @(\"907\u00f73=\", \"974\u00f77=\"),\n    @(\"571\u00f79=\", \"410\u00f72=\"),\n    @(\"235\u00f76=\", \"473\u00f78=\"),\n    @(\"684\u00f76=\", \"929\u00f73=\"),\n    @(\"449\u00f74=\", \"708\u00f77=\"),\n    @(\"598\u00f75=\", \"790\u00f76=\"),\n    @(\"758\u00f77=\", \"360\u00f76=\"),\n    @(\"109\u00f74=\", \"855\u00f77=\"),\n    @(\"361\u00f76=\", \"958\u00f79=\"),\n    @(\"276\u00f74=\", \"143\u00f76=\"),\n    @(\"986\u00f76=\", \"565\u00f74=\"),\n    @(\"696\u00f78=\", \"678\u00f77=\"),\n    @(\"623\u00f72=\", \"812\u00f75=\"),\n    @(\"254\u00f77=\", \"193\u00f72=\"),\n    @(\"740\u00f79=\", \"924\u00f76=\"),\n    @(\"500\u00f74=\", \"888\u00f72=\"),\n    @(\"950\u00f78=\", \"881\u00f72=\"),\n    @(\"838\u00f76=\", \"875\u00f77=\"),\n    @(\"865\u00f77=\", \"740\u00f74=\"),\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $d.Content.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}"}
